# Commit: "Tue, Mar 31, 2020  7:04:58 AM"
#
# The deck's Slide Master theme ("Integral" / "Red Violet" colour scheme,
# persisted as ppt/theme/theme2.xml) is recoloured to the standard
# "Office" colour scheme (the palette that, in the source deck, lives in
# the otherwise-unused ppt/theme/theme1.xml referenced only by the Notes
# Master).
#
# PowerPoint's ColorScheme.Colors(n).RGB setter takes a COM "OLE_COLOR"
# (0xBBGGRR, i.e. the RGB bytes reversed), so we convert each target hex
# colour (RRGGBB) into that packed integer before assigning it.

function ToOleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - in clrScheme document order
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $cs.Colors($i).RGB = ToOleColor $officeColors[$i - 1]
}
